$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 is the last existing data row; append two more rows (16, 17) that
# duplicate its values, matching the target diff.
$values = @(
    "area undefined",
    5,
    100,
    100,
    124517501.1612876,
    100000000,
    24517501.16128756,
    16332027.26010531,
    116302.8620469598,
    2956469.651481002,
    1049740.5,
    11539361.24657734,
    615203,
    54950,
    8185473.901182257,
    19.69000416216997,
    12191684.05266326,
    49.72645447209597,
    12307986.91471022,
    50.20082117562691,
    12316026.91471022,
    50.2336140771021
)

foreach ($r in 16..17) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
}
